$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the manual annotation scores (columns E:J, rows 2-22) ---
$data = New-Object 'object[,]' 21,6
$data[0,0] = 2
$data[0,1] = 2
$data[0,2] = 2
$data[0,3] = 2
$data[0,4] = 2
$data[0,5] = 2
$data[1,0] = 2
$data[1,1] = 1
$data[1,2] = 1
$data[1,3] = 1
$data[1,4] = 1
$data[1,5] = 2
$data[2,0] = 2
$data[2,1] = 2
$data[2,2] = 2
$data[2,3] = 2
$data[2,4] = 2
$data[2,5] = 2
$data[3,0] = 2
$data[3,1] = 1
$data[3,2] = 1
$data[3,3] = 1
$data[3,4] = 1
$data[3,5] = 2
$data[4,0] = 2
$data[4,1] = 2
$data[4,2] = 2
$data[4,3] = 2
$data[4,4] = 2
$data[4,5] = 2
$data[5,0] = 2
$data[5,1] = 1
$data[5,2] = 2
$data[5,3] = 2
$data[5,4] = 1
$data[5,5] = 2
$data[6,0] = 2
$data[6,1] = 2
$data[6,2] = 1
$data[6,3] = 1
$data[6,4] = 1
$data[6,5] = 2
$data[7,0] = 2
$data[7,1] = 2
$data[7,2] = 2
$data[7,3] = 2
$data[7,4] = 2
$data[7,5] = 2
$data[8,0] = 2
$data[8,1] = 2
$data[8,2] = 1
$data[8,3] = 2
$data[8,4] = 2
$data[8,5] = 2
$data[9,0] = 2
$data[9,1] = 1
$data[9,2] = 2
$data[9,3] = 1
$data[9,4] = 2
$data[9,5] = 2
$data[10,0] = 2
$data[10,1] = 2
$data[10,2] = 2
$data[10,3] = 2
$data[10,4] = 2
$data[10,5] = 2
$data[11,0] = 2
$data[11,1] = 2
$data[11,2] = 2
$data[11,3] = 2
$data[11,4] = 2
$data[11,5] = 2
$data[12,0] = 2
$data[12,1] = 1
$data[12,2] = 1
$data[12,3] = 1
$data[12,4] = 1
$data[12,5] = 2
$data[13,0] = 2
$data[13,1] = 2
$data[13,2] = 2
$data[13,3] = 2
$data[13,4] = 2
$data[13,5] = 2
$data[14,0] = 2
$data[14,1] = 1
$data[14,2] = 2
$data[14,3] = 1
$data[14,4] = 1
$data[14,5] = 2
$data[15,0] = 2
$data[15,1] = 2
$data[15,2] = 1
$data[15,3] = 1
$data[15,4] = 2
$data[15,5] = 2
$data[16,0] = 2
$data[16,1] = 2
$data[16,2] = 1
$data[16,3] = 1
$data[16,4] = 2
$data[16,5] = 2
$data[17,0] = 2
$data[17,1] = 2
$data[17,2] = 1
$data[17,3] = 1
$data[17,4] = 2
$data[17,5] = 2
$data[18,0] = 2
$data[18,1] = 2
$data[18,2] = 2
$data[18,3] = 2
$data[18,4] = 2
$data[18,5] = 2
$data[19,0] = 2
$data[19,1] = 2
$data[19,2] = 2
$data[19,3] = 2
$data[19,4] = 2
$data[19,5] = 2
$data[20,0] = 2
$data[20,1] = 2
$data[20,2] = 2
$data[20,3] = 2
$data[20,4] = 2
$data[20,5] = 2

$ws.Range("E2:J22").Value = $data

# --- Window / view adjustments ---
$ws.Activate()

# Scroll / freeze top row, with the view left on row 21 and E20 selected
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E20").Select()

# Zoom to 85%
$excel.ActiveWindow.Zoom = 85

Write-Output "done"
